$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 121. This shifts all existing rows 121-186
# down to 122-187 (matching the diff, where every row's data is simply
# "pushed down" by one position and a brand new record appears at the top
# of the block, row 121).
$ws.Rows.Item(121).Insert()

# Populate the newly inserted row 121 with the new weekly record.
$ws.Cells.Item(121, 1).Value = 5
$ws.Cells.Item(121, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(121, 3).Value = "Maule"
$ws.Cells.Item(121, 4).Value = 44488
$ws.Cells.Item(121, 5).Value = 7
$ws.Cells.Item(121, 6).Value = 100112009
$ws.Cells.Item(121, 7).Value = "Acelga"
$ws.Cells.Item(121, 8).Value = "Sin especificar"
$ws.Cells.Item(121, 9).Value = "Primera"
$ws.Cells.Item(121, 10).Value = 500
$ws.Cells.Item(121, 11).Value = 2000
$ws.Cells.Item(121, 12).Value = 2000
$ws.Cells.Item(121, 13).Value = 2000
$ws.Cells.Item(121, 14).Value = "`$/docena de atados (4 kilos)"
$ws.Cells.Item(121, 15).Value = "Región del Maule"
$ws.Cells.Item(121, 16).Value = 500
$ws.Cells.Item(121, 17).Value = 4
$ws.Cells.Item(121, 18).Value = "Hortaliza"
